# Weekly update: insert a new price record as row 6 (shifting existing
# rows 6-47 down to 7-48) on the single worksheet of the "Poroto granado"
# subset workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6, pushing the rest of the table down.
$ws.Rows(6).Insert()

# Populate the new row 6 with the latest weekly record.
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44552
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112030
$ws.Range("G6").Value = "Poroto granado"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 35000
$ws.Range("L6").Value = 36000
$ws.Range("M6").Value = 35500
$ws.Range("N6").Value = "`$/malla 25 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 1420
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
